# overzicht-inzet.xlsx - "DDM1 fragments volgens railwiki"
#
# Adds a breakdown table (rows 37-43) under the existing "Consists" header
# (A36) on the "DDM1" sheet, listing each DDM1 consist-length "fragment"
# with its first/last date, plus a note-comment on B40 explaining that no
# data could be found for earlier consists.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DDM1")

# --- header row (A36 already holds "Consists") : add column headers ---
$ws.Range("B36").Value = "Eerste"
$ws.Range("C36").Value = "Laatste"

# --- new data rows 37-43 -------------------------------------------------
# Label (column A) cells are written in the same order the original author
# used so new shared-string entries land in the same slots as the source
# workbook (55=DDM1 7 1985 ... 61=DDM1 4 2016).
$ws.Range("A41").Value = "DDM1 7 1985"
$ws.Range("A39").Value = "DDM1 6 1985"
$ws.Range("A42").Value = "DDM1 7 1993"
$ws.Range("A38").Value = "DDM1 5 1995"
$ws.Range("A43").Value = "DDM1 7 1999"
$ws.Range("A40").Value = "DDM1 6 2004"
$ws.Range("A37").Value = "DDM1 4 2016"

# Eerste / Laatste date pairs (column B / C), formatted like the existing
# date columns elsewhere on this sheet (numFmt "mmm-yy" -> style index 2).
$dates = @{
    37 = @(42522, 43800)
    38 = @(34851, 37591)
    39 = @(31199, 34851)
    40 = @(37043, 41244)
    41 = @(31199, 34121)
    42 = @(34121, 36312)
    43 = @(36312, 37956)
}

foreach ($row in @(37, 38, 39, 40, 41, 42, 43)) {
    $pair = $dates[$row]

    $bCell = $ws.Range("B$row")
    $bCell.Value = $pair[0]
    $bCell.NumberFormat = "mmm-yy"

    $cCell = $ws.Range("C$row")
    $cCell.Value = $pair[1]
    $cCell.NumberFormat = "mmm-yy"
}

# --- comment on B40 -------------------------------------------------------
$ws.Range("B40").AddComment("Tom:`nkan geen data vinden over eerdere stammen")

# --- selection / view state ----------------------------------------------
$ws.Activate()
$ws.Range("E37").Select()
